$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "147"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "394886.40"

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "422"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1097510.82"

$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "149"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "339041.00"

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "797"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3004058.81"

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "169"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "527316.18"

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "94"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "230800.00"

$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "94"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "229788.98"

$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "186"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "421089.87"

$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "14"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "54000.00"

$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "35"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "101891.77"

$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "129"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "366175.00"

$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "14"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "35800.00"

$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "312"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1285511.12"

$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "5"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13000.00"

$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "13"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "39500.00"

$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "7"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "23500.00"

$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "45"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "213018.00"

$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "55"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "137310.00"

$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "168"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "475408.00"

$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "377"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1499848.18"

$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "26"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "78500.00"

$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "76"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "171893.00"

$ws.Range("C74").NumberFormat = "@"
$ws.Range("C74").Value = "20"
$ws.Range("D74").NumberFormat = "@"
$ws.Range("D74").Value = "79400.00"

$ws.Range("C77").NumberFormat = "@"
$ws.Range("C77").Value = "90"
$ws.Range("D77").NumberFormat = "@"
$ws.Range("D77").Value = "238487.00"

$ws.Range("C78").NumberFormat = "@"
$ws.Range("C78").Value = "211"
$ws.Range("D78").NumberFormat = "@"
$ws.Range("D78").Value = "587693.00"

$ws.Range("C80").NumberFormat = "@"
$ws.Range("C80").Value = "491"
$ws.Range("D80").NumberFormat = "@"
$ws.Range("D80").Value = "2147734.03"

$ws.Range("C81").NumberFormat = "@"
$ws.Range("C81").Value = "14"
$ws.Range("D81").NumberFormat = "@"
$ws.Range("D81").Value = "39000.00"

$ws.Range("C82").NumberFormat = "@"
$ws.Range("C82").Value = "8"
$ws.Range("D82").NumberFormat = "@"
$ws.Range("D82").Value = "23500.00"

$ws.Range("C83").NumberFormat = "@"
$ws.Range("C83").Value = "18"
$ws.Range("D83").NumberFormat = "@"
$ws.Range("D83").Value = "58069.00"

$ws.Range("C84").NumberFormat = "@"
$ws.Range("C84").Value = "73"
$ws.Range("D84").NumberFormat = "@"
$ws.Range("D84").Value = "255072.36"

$ws.Range("C85").NumberFormat = "@"
$ws.Range("C85").Value = "36"
$ws.Range("D85").NumberFormat = "@"
$ws.Range("D85").Value = "116669.00"

$ws.Range("C88").NumberFormat = "@"
$ws.Range("C88").Value = "73"
$ws.Range("D88").NumberFormat = "@"
$ws.Range("D88").Value = "333636.08"

$ws.Range("C123").NumberFormat = "@"
$ws.Range("C123").Value = "116"
$ws.Range("D123").NumberFormat = "@"
$ws.Range("D123").Value = "304781.45"

$ws.Range("C124").NumberFormat = "@"
$ws.Range("C124").Value = "496"
$ws.Range("D124").NumberFormat = "@"
$ws.Range("D124").Value = "2224983.06"

$ws.Range("C125").NumberFormat = "@"
$ws.Range("C125").Value = "8"
$ws.Range("D125").NumberFormat = "@"
$ws.Range("D125").Value = "27000.00"

$ws.Range("C127").NumberFormat = "@"
$ws.Range("C127").Value = "32"
$ws.Range("D127").NumberFormat = "@"
$ws.Range("D127").Value = "74500.00"

$ws.Range("C129").NumberFormat = "@"
$ws.Range("C129").Value = "44"
$ws.Range("D129").NumberFormat = "@"
$ws.Range("D129").Value = "166579.76"

$ws.Range("C130").NumberFormat = "@"
$ws.Range("C130").Value = "58"
$ws.Range("D130").NumberFormat = "@"
$ws.Range("D130").Value = "165793.82"

$ws.Range("C133").NumberFormat = "@"
$ws.Range("C133").Value = "127"
$ws.Range("D133").NumberFormat = "@"
$ws.Range("D133").Value = "328364.44"

$ws.Range("C201").NumberFormat = "@"
$ws.Range("C201").Value = "662"
$ws.Range("D201").NumberFormat = "@"
$ws.Range("D201").Value = "2525282.58"
